$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2: Bitcoin
Set-TextCell 2 2 'Bitcoin'
Set-TextCell 2 3 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextCell 2 4 '26.486.20'
Set-TextCell 2 5 '  -0.06%  '

# Row 3: Ethereum
Set-TextCell 3 2 'Ethereum'
Set-TextCell 3 3 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextCell 3 4 '1.725.67'
Set-TextCell 3 5 '  -0.50%  '

# Row 4: TetherUSD
Set-TextCell 4 2 'TetherUSD'
Set-TextCell 4 3 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextCell 4 4 '0.9972'
Set-TextCell 4 5 '  -0.30%  '

# Row 5: BNB
Set-TextCell 5 2 'BNB'
Set-TextCell 5 3 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell 5 4 '242.88'
Set-TextCell 5 5 '  -1.71%  '

# Row 6: USDC
Set-TextCell 6 2 'USDC'
Set-TextCell 6 3 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextCell 6 4 '0.9976'
Set-TextCell 6 5 '  -0.29%  '

# Row 7: XRP
Set-TextCell 7 2 'XRP'
Set-TextCell 7 3 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextCell 7 4 '0.4901'
Set-TextCell 7 5 '  +0.36%  '

# Row 8: Cardano
Set-TextCell 8 2 'Cardano'
Set-TextCell 8 3 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell 8 4 '0.2609'
Set-TextCell 8 5 '  -2.07%  '

# Row 9: Dogecoin
Set-TextCell 9 2 'Dogecoin'
Set-TextCell 9 3 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 9 4 '0.06200'
Set-TextCell 9 5 '  -0.25%  '

# Row 10: WrappedEther
Set-TextCell 10 2 'WrappedEther'
Set-TextCell 10 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 10 4 '1.728.43'
Set-TextCell 10 5 '  -0.06%  '

# Row 11: TRON
Set-TextCell 11 2 'TRON'
Set-TextCell 11 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 11 4 '0.06967'
Set-TextCell 11 5 '  -0.96%  '

# Row 12: Solana
Set-TextCell 12 2 'Solana'
Set-TextCell 12 3 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 12 4 '15.64'
Set-TextCell 12 5 '  -0.11%  '

# Row 13: Polkadot
Set-TextCell 13 2 'Polkadot'
Set-TextCell 13 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 13 4 '4.529'
Set-TextCell 13 5 '  -1.30%  '

# Row 14: Polygon
Set-TextCell 14 2 'Polygon'
Set-TextCell 14 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 14 4 '0.6015'
Set-TextCell 14 5 '  -1.09%  '

# Row 15: Litecoin
Set-TextCell 15 2 'Litecoin'
Set-TextCell 15 3 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 15 4 '77.17'
Set-TextCell 15 5 '  -0.36%  '

# Row 16: Dai
Set-TextCell 16 2 'Dai'
Set-TextCell 16 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 16 4 '0.9972'
Set-TextCell 16 5 '  -0.26%  '

# Row 17: WrappedBTC
Set-TextCell 17 2 'WrappedBTC'
Set-TextCell 17 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 17 4 '26.478.83'
Set-TextCell 17 5 '  -0.10%  '

# Row 18: BinanceUSD
Set-TextCell 18 2 'BinanceUSD'
Set-TextCell 18 3 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 18 4 '0.9971'
Set-TextCell 18 5 '  -0.33%  '

# Row 19: ShibaInu
Set-TextCell 19 2 'ShibaInu'
Set-TextCell 19 3 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 19 4 '0.000007155'
Set-TextCell 19 5 '  -2.67%  '

# Row 20: Avalanche
Set-TextCell 20 2 'Avalanche'
Set-TextCell 20 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 20 4 '11.33'
Set-TextCell 20 5 '  -1.59%  '

# Row 21: WrappedliquidstakedEther2.0
Set-TextCell 21 2 'WrappedliquidstakedEther2.0'
Set-TextCell 21 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 21 4 '1.946.37'
Set-TextCell 21 5 '  -0.37%  '

# Row 22: Uniswap
Set-TextCell 22 2 'Uniswap'
Set-TextCell 22 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 22 4 '4.452'
Set-TextCell 22 5 '  -2.26%  '

# Row 23: Cosmos
Set-TextCell 23 2 'Cosmos'
Set-TextCell 23 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 23 4 '8.522'
Set-TextCell 23 5 '  -2.55%  '

# Row 24: Chainlink
Set-TextCell 24 2 'Chainlink'
Set-TextCell 24 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 24 4 '5.112'
Set-TextCell 24 5 '  -2.23%  '

# Row 25: Monero
Set-TextCell 25 2 'Monero'
Set-TextCell 25 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 25 4 '137.64'
Set-TextCell 25 5 '  -2.31%  '

# Row 26: EthereumClassic
Set-TextCell 26 2 'EthereumClassic'
Set-TextCell 26 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 26 4 '15.29'
Set-TextCell 26 5 '  -1.11%  '

# Row 27: Toncoin
Set-TextCell 27 2 'Toncoin'
Set-TextCell 27 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 27 4 '1.410'
Set-TextCell 27 5 '  -0.06%  '

# Row 28: LidoDAOToken
Set-TextCell 28 2 'LidoDAOToken'
Set-TextCell 28 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 28 4 '1.750'
Set-TextCell 28 5 '  -1.21%  '

# Row 29: BitcoinCash
Set-TextCell 29 2 'BitcoinCash'
Set-TextCell 29 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 29 4 '106.56'
Set-TextCell 29 5 '  -1.44%  '

# Row 30: InternetComputer(DFINITY)
Set-TextCell 30 2 'InternetComputer(DFINITY)'
Set-TextCell 30 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 30 4 '3.923'
Set-TextCell 30 5 '  -2.29%  '

# Row 31: Stellar
Set-TextCell 31 2 'Stellar'
Set-TextCell 31 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 31 4 '0.08002'
Set-TextCell 31 5 '  -0.51%  '

# Row 32: Filecoin
Set-TextCell 32 2 'Filecoin'
Set-TextCell 32 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 32 4 '3.640'
Set-TextCell 32 5 '  -1.27%  '

# Row 33: Hedera
Set-TextCell 33 2 'Hedera'
Set-TextCell 33 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 33 4 '0.04491'
Set-TextCell 33 5 '  -1.57%  '

# Row 34: HuobiToken
Set-TextCell 34 2 'HuobiToken'
Set-TextCell 34 3 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 34 4 '2.598'
Set-TextCell 34 5 '  -0.58%  '

# Row 35: ARBITRUM
Set-TextCell 35 2 'ARBITRUM'
Set-TextCell 35 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 35 4 '1.003'
Set-TextCell 35 5 '  -0.32%  '

# Row 36: ImmutableX
Set-TextCell 36 2 'ImmutableX'
Set-TextCell 36 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 36 4 '0.6229'
Set-TextCell 36 5 '  -1.66%  '

# Row 37: TrustWalletToken
Set-TextCell 37 2 'TrustWalletToken'
Set-TextCell 37 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 37 4 '0.9341'
Set-TextCell 37 5 '  +4.12%  '

# Row 38: RenderToken
Set-TextCell 38 2 'RenderToken'
Set-TextCell 38 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 38 4 '1.995'
Set-TextCell 38 5 '  -1.28%  '

# Row 39: MXToken
Set-TextCell 39 2 'MXToken'
Set-TextCell 39 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 39 4 '2.386'
Set-TextCell 39 5 '  -0.45%  '

# Row 40: PaxDollar
Set-TextCell 40 2 'PaxDollar'
Set-TextCell 40 3 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 40 4 '0.9969'
Set-TextCell 40 5 '  -0.76%  '

# Row 41: VeChain
Set-TextCell 41 2 'VeChain'
Set-TextCell 41 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 41 4 '0.01483'
Set-TextCell 41 5 '  -1.27%  '

# Row 42: Quant
Set-TextCell 42 2 'Quant'
Set-TextCell 42 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 42 4 '99.82'
Set-TextCell 42 5 '  -1.89%  '

# Row 43: FraxShare
Set-TextCell 43 2 'FraxShare'
Set-TextCell 43 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 43 4 '5.394'
Set-TextCell 43 5 '  -0.34%  '

# Row 44: TheSandbox
Set-TextCell 44 2 'TheSandbox'
Set-TextCell 44 3 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 44 4 '0.3848'
Set-TextCell 44 5 '  -0.97%  '

# Row 45: Aptos
Set-TextCell 45 2 'Aptos'
Set-TextCell 45 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 45 4 '6.898'
Set-TextCell 45 5 '  -0.35%  '

# Row 46: Algorand
Set-TextCell 46 2 'Algorand'
Set-TextCell 46 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 46 4 '0.1159'
Set-TextCell 46 5 '  -2.08%  '

# Row 47: Cronos
Set-TextCell 47 2 'Cronos'
Set-TextCell 47 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 47 4 '0.05369'
Set-TextCell 47 5 '  -0.51%  '

# Row 48: Elrond
Set-TextCell 48 2 'Elrond'
Set-TextCell 48 3 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell 48 4 '30.56'
Set-TextCell 48 5 '  +0.30%  '

# Row 49: EnergySwap
Set-TextCell 49 2 'EnergySwap'
Set-TextCell 49 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 49 4 '7.741'
Set-TextCell 49 5 '  -0.96%  '

# Row 50: Aave
Set-TextCell 50 2 'Aave'
Set-TextCell 50 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 50 4 '51.43'
Set-TextCell 50 5 '  -0.39%  '

# Row 51: NEARProtocol
Set-TextCell 51 2 'NEARProtocol'
Set-TextCell 51 3 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 51 4 '1.226'
Set-TextCell 51 5 '  -2.27%  '
